$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "copyright" column entirely (column A), shifting B:H left to A:G
$ws.Range("A1:A2").Delete(-4159) | Out-Null  # xlShiftToLeft

# Update header row (row 1) with the new values
$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "explanation"
$ws.Range("C1").Value = "hdurl"
$ws.Range("D1").Value = "media_type"
$ws.Range("E1").Value = "service_version"
$ws.Range("F1").Value = "title"
$ws.Range("G1").Value = "url"

# Update data row (row 2) with the new values
# A2 holds a date-like string; force text storage so it is not
# auto-converted into a numeric date serial value.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-03-01"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = "With spacecraft thrusters at top center, the rugged surface of the Moon lies below the Blue Ghost lander in this space age video frame. The view of the lunar far side was captured by the Firefly Aerospace lunar lander on February 24, following a maneuver to circularize its orbit about 100 kilometers above the lunar surface. The robotic lunar lander is scheduled to touch down tomorrow, Sunday, March 2, at 3:34am Eastern Time in the Mare Crisium impact basin on the lunar near side. In support of the Artemis campaign, Blue Ghost is set to deliver science and technology experiments to the Moon, part of NASA's Commercial Lunar Payload Services program. Blue Ghost's mission on the surface is planned to operate during the lunar daylight hours at the landing site, about 14 Earth days."
$ws.Range("C2").Value = "https://apod.nasa.gov/apod/image/2503/BlueGhost_lunar1067.jpg"
$ws.Range("D2").Value = "image"
$ws.Range("E2").Value = "v1"
$ws.Range("F2").Value = "Blue Ghost to the Moon"
$ws.Range("G2").Value = "https://apod.nasa.gov/apod/image/2503/BlueGhost_lunar1067.jpg"
